# Auto-update draw results: append the 2025-12-15 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

# B/D/E are plain, non-numeric-looking text, so a direct .Value assignment
# is stored as text without Excel re-interpreting it.
$ws.Range("B$row").Value = "Pick 3"
$ws.Range("D$row").Value = "7-3-8"
$ws.Range("E$row").Value = "2025-12-15T21:46:52.959+04:00"

# A ("2025-12-15") and C ("251215") look like a date / a number, so a plain
# .Value assignment would get auto-coerced into a date serial / numeric
# value. Write them as text-literal formulas first (forces text), then
# copy/paste-special as values so the cells end up holding plain text with
# the sheet's default (unformatted) style - exactly like the rest of the
# column - instead of a cached formula or a cell carrying an explicit
# "@" / quote-prefix style.
$ws.Range("A$row").Formula = "=""2025-12-15"""
$ws.Range("C$row").Formula = "=""251215"""

$rng = $ws.Range("A$($row):E$row")
$rng.Copy()
$rng.PasteSpecial(-4163)
